# Updates betting-odds cell values on Sheet1 for the FlashScore 2024-10-09
# "Jogos da Semana" workbook (rows 8-11), matching the upstream commit
# "Atualizando o arquivo XLSX".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("G8").Value = 2.15
$ws.Range("H8").Value = 3.4
$ws.Range("L8").Value = 3.7
$ws.Range("N8").Value = 7.8
$ws.Range("O8").Value = 1.25
$ws.Range("P8").Value = 3.6
$ws.Range("Q8").Value = 1.75
$ws.Range("R8").Value = 2
$ws.Range("S8").Value = 1.39
$ws.Range("U8").Value = 1.62
$ws.Range("V8").Value = 2.18
$ws.Range("AB8").Value = 23
$ws.Range("AC8").Value = 7.8
$ws.Range("AD8").Value = 6.6
$ws.Range("AG8").Value = 350
$ws.Range("AH8").Value = 10.75
$ws.Range("AI8").Value = 17
$ws.Range("AJ8").Value = 11
$ws.Range("AL8").Value = 26
$ws.Range("AM8").Value = 30
$ws.Range("AP8").Value = 18
$ws.Range("AU8").Value = 6.8
$ws.Range("AV8").Value = 55
$ws.Range("AY8").Value = 24
$ws.Range("AZ8").Value = 90
$ws.Range("BA8").Value = 120

# Row 9
$ws.Range("H9").Value = 3.5
$ws.Range("I9").Value = 3.4
$ws.Range("L9").Value = 3.8
$ws.Range("P9").Value = 3.65
$ws.Range("S9").Value = 1.37
$ws.Range("AH9").Value = 12.5
$ws.Range("AU9").Value = 6.8

# Row 10
$ws.Range("G10").Value = 2.62
$ws.Range("H10").Value = 3.15
$ws.Range("J10").Value = 3.25
$ws.Range("K10").Value = 2.05
$ws.Range("N10").Value = 7.1
$ws.Range("O10").Value = 1.3
$ws.Range("P10").Value = 3.2
$ws.Range("Q10").Value = 1.9
$ws.Range("R10").Value = 1.83
$ws.Range("T10").Value = 2.65
$ws.Range("W10").Value = 8.5
$ws.Range("Z10").Value = 30
$ws.Range("AA10").Value = 22
$ws.Range("AB10").Value = 29
$ws.Range("AC10").Value = 7.1
$ws.Range("AD10").Value = 6.1
$ws.Range("AH10").Value = 9
$ws.Range("AI10").Value = 14
$ws.Range("AM10").Value = 27
$ws.Range("AN10").Value = 4.6
$ws.Range("AO10").Value = 14.5
$ws.Range("AP10").Value = 22
$ws.Range("AQ10").Value = 65
$ws.Range("AR10").Value = 100
$ws.Range("AS10").Value = 300
$ws.Range("AT10").Value = 2.65
$ws.Range("AW10").Value = 4.55
$ws.Range("AX10").Value = 14
$ws.Range("AY10").Value = 21
$ws.Range("AZ10").Value = 60
$ws.Range("BA10").Value = 90
$ws.Range("BB10").Value = 300

# Row 11
$ws.Range("H11").Value = 4.3
$ws.Range("I11").Value = 6.9
$ws.Range("K11").Value = 2.27
$ws.Range("L11").Value = 6.5
$ws.Range("N11").Value = 7.7
$ws.Range("Y11").Value = 8.5
$ws.Range("AB11").Value = 32
$ws.Range("AC11").Value = 7.7
$ws.Range("AD11").Value = 8.5
$ws.Range("AH11").Value = 16.5
$ws.Range("AI11").Value = 40
$ws.Range("AJ11").Value = 22
$ws.Range("AK11").Value = 150
$ws.Range("AL11").Value = 80
$ws.Range("AV11").Value = 100
$ws.Range("AW11").Value = 7.9
$ws.Range("AX11").Value = 40
$ws.Range("AZ11").Value = 300
